$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1644.36
$ws.Range("I33").Value = 979.2105
$ws.Range("K33").Value = 979.2105
$ws.Range("M33").Value = -750.2105
$ws.Range("H40").Value = 1044.1666
$ws.Range("I40").Value = 1040
$ws.Range("J40").Value = 1056.6666
$ws.Range("K40").Value = 1040
$ws.Range("L40").Value = 1056.6666
$ws.Range("M40").Value = -865
$ws.Range("N40").Value = -1406.6666
$ws.Range("H70").Value = 28146.65
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 31162.945
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 93488.83499999999
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -94028.83499999999
$ws.Range("H73").Value = 28146.65
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 31162.945
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 93488.83499999999
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -95360.83499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 2875
$ws.Range("J27").Value = 2875
$ws.Range("L27").Value = 2875
$ws.Range("N27").Value = -3243
$ws.Range("H32").Value = 4414.11
$ws.Range("I32").Value = 4241.7324
$ws.Range("K32").Value = 4241.7324
$ws.Range("M32").Value = -3954.7324
$ws.Range("H61").Value = 4709248.5
$ws.Range("I61").Value = 3352006.5
$ws.Range("J61").Value = 11766906
$ws.Range("K61").Value = 3352006.5
$ws.Range("L61").Value = 11766906
$ws.Range("M61").Value = -3351794.5
$ws.Range("N61").Value = -11767330
$ws.Range("H74").Value = 45409930
$ws.Range("I74").Value = 50284224
$ws.Range("J74").Value = 38098484
$ws.Range("K74").Value = 50284224
$ws.Range("L74").Value = 38098484
$ws.Range("M74").Value = -50283350
$ws.Range("N74").Value = -38100232
$ws.Range("H77").Value = 45409930
$ws.Range("I77").Value = 50284224
$ws.Range("J77").Value = 38098484
$ws.Range("K77").Value = 251421120
$ws.Range("L77").Value = 190492420
$ws.Range("M77").Value = -251416752
$ws.Range("N77").Value = -190501156
$ws.Range("H136").Value = 4709248.5
$ws.Range("I136").Value = 3352006.5
$ws.Range("J136").Value = 11766906
$ws.Range("K136").Value = 10056019.5
$ws.Range("L136").Value = 35300718
$ws.Range("M136").Value = -10053469.5
$ws.Range("N136").Value = -35305818

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 7500
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10228
$ws.Range("H134").Value = 18256588
$ws.Range("I134").Value = 20000970
$ws.Range("J134").Value = 7354194.5
$ws.Range("K134").Value = 60002910
$ws.Range("L134").Value = 22062583.5
$ws.Range("M134").Value = -60000375
$ws.Range("N134").Value = -22067653.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = 0
$ws.Range("H31").Value = 14927052
$ws.Range("I31").Value = 35715652
$ws.Range("J31").Value = 1902.6666
$ws.Range("K31").Value = 35715652
$ws.Range("L31").Value = 1902.6666
$ws.Range("M31").Value = -35715357
$ws.Range("N31").Value = -2492.6666
$ws.Range("H34").Value = 14927052
$ws.Range("I34").Value = 35715652
$ws.Range("J34").Value = 1902.6666
$ws.Range("K34").Value = 35715652
$ws.Range("L34").Value = 1902.6666
$ws.Range("M34").Value = -35715450
$ws.Range("N34").Value = -2306.6666
$ws.Range("H36").Value = 18900.143
$ws.Range("I36").Value = 4764.25
$ws.Range("J36").Value = 37748
$ws.Range("K36").Value = 4764.25
$ws.Range("L36").Value = 37748
$ws.Range("M36").Value = -4376.25
$ws.Range("N36").Value = -38524
$ws.Range("H40").Value = 18900.143
$ws.Range("I40").Value = 4764.25
$ws.Range("J40").Value = 37748
$ws.Range("K40").Value = 4764.25
$ws.Range("L40").Value = 37748
$ws.Range("M40").Value = -4604.25
$ws.Range("N40").Value = -38068
$ws.Range("H58").Value = 2976924
$ws.Range("I58").Value = 3572178.8
$ws.Range("J58").Value = 650
$ws.Range("K58").Value = 3572178.8
$ws.Range("L58").Value = 650
$ws.Range("M58").Value = -3571975.8
$ws.Range("N58").Value = -1056
$ws.Range("H105").Value = 3886.7026
$ws.Range("I105").Value = 923.38464
$ws.Range("J105").Value = 10890.909
$ws.Range("K105").Value = 923.38464
$ws.Range("L105").Value = 10890.909
$ws.Range("M105").Value = 823.61536
$ws.Range("N105").Value = -14384.909
$ws.Range("H120").Value = 86209.664
$ws.Range("J120").Value = 29166.5
$ws.Range("L120").Value = 29166.5
$ws.Range("N120").Value = -36424.5
$ws.Range("H136").Value = 2976924
$ws.Range("I136").Value = 3572178.8
$ws.Range("J136").Value = 650
$ws.Range("K136").Value = 10716536.4
$ws.Range("L136").Value = 1950
$ws.Range("M136").Value = -10713986.4
$ws.Range("N136").Value = -7050
$ws.Range("H141").Value = 225179.22
$ws.Range("J141").Value = 252451.62
$ws.Range("L141").Value = 252451.62
$ws.Range("N141").Value = -262811.62
$ws.Range("N23").Value = ""
$ws.Range("N27").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2565189.2
$ws.Range("J5").Value = 3032271.5
$ws.Range("L5").Value = 9096814.5
$ws.Range("N5").Value = -9097038.5
$ws.Range("H68").Value = 4161.1763
$ws.Range("I68").Value = 680
$ws.Range("J68").Value = 9784.615
$ws.Range("K68").Value = 2040
$ws.Range("L68").Value = 29353.845
$ws.Range("M68").Value = -1229
$ws.Range("N68").Value = -30975.845
$ws.Range("H71").Value = 4161.1763
$ws.Range("I71").Value = 680
$ws.Range("J71").Value = 9784.615
$ws.Range("K71").Value = 6120
$ws.Range("L71").Value = 88061.535
$ws.Range("M71").Value = -2064
$ws.Range("N71").Value = -96173.535
$ws.Range("H113").Value = 2545.8508
$ws.Range("I113").Value = 3086.282
$ws.Range("K113").Value = 9258.846000000001
$ws.Range("M113").Value = -7088.846000000001
$ws.Range("H135").Value = 2565189.2
$ws.Range("J135").Value = 3032271.5
$ws.Range("L135").Value = 27290443.5
$ws.Range("N135").Value = -27295513.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4314226
$ws.Range("I70").Value = 2142296.2
$ws.Range("K70").Value = 2142296.2
$ws.Range("M70").Value = -2142026.2
$ws.Range("H73").Value = 4314226
$ws.Range("I73").Value = 2142296.2
$ws.Range("K73").Value = 2142296.2
$ws.Range("M73").Value = -2141360.2
$ws.Range("H122").Value = 5976089
$ws.Range("I122").Value = 29395.908
$ws.Range("J122").Value = 27780630
$ws.Range("K122").Value = 88187.724
$ws.Range("L122").Value = 83341890
$ws.Range("M122").Value = -85737.724
$ws.Range("N122").Value = -83346790
$ws.Range("H132").Value = 11043449
$ws.Range("I132").Value = 7188270
$ws.Range("K132").Value = 21564810
$ws.Range("M132").Value = -21562280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 6610.25
$ws.Range("I34").Value = 2501
$ws.Range("K34").Value = 2501
$ws.Range("M34").Value = -2329
$ws.Range("H68").Value = 2633.5334
$ws.Range("I68").Value = 2166.6667
$ws.Range("J68").Value = 2750.25
$ws.Range("K68").Value = 2166.6667
$ws.Range("L68").Value = 2750.25
$ws.Range("M68").Value = -1417.6667
$ws.Range("N68").Value = -4248.25
$ws.Range("H71").Value = 2633.5334
$ws.Range("I71").Value = 2166.6667
$ws.Range("J71").Value = 2750.25
$ws.Range("K71").Value = 10833.3335
$ws.Range("L71").Value = 13751.25
$ws.Range("M71").Value = -7089.333500000001
$ws.Range("N71").Value = -21239.25
$ws.Range("H132").Value = 2409218
$ws.Range("I132").Value = 3082562
$ws.Range("J132").Value = 4418
$ws.Range("K132").Value = 9247686
$ws.Range("L132").Value = 13254
$ws.Range("M132").Value = -9245156
$ws.Range("N132").Value = -18314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").Value = 0
$ws.Range("H136").Value = 6920.676
$ws.Range("I136").Value = 4753.5557
$ws.Range("J136").Value = 12771.9
$ws.Range("K136").Value = 14260.6671
$ws.Range("L136").Value = 38315.7
$ws.Range("M136").Value = -11710.6671
$ws.Range("N136").Value = -43415.7
$ws.Range("N40").Value = ""
